# Negative scenarios for user details
# Adds a new worksheet "PutUserIdNegative" (copied/trimmed from "PutUser")
# containing two negative test rows for invalid userId / invalid timezone,
# positioned right after "PutProgramBatchStatusNegative".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet by copying "PutUser" (it already has the right
#    column styles/hyperlinks for most of the columns we need) and place
#    it right after "PutProgramBatchStatusNegative".
# ---------------------------------------------------------------------
$srcSheet   = $wb.Worksheets.Item("PutUser")
$afterSheet = $wb.Worksheets.Item("PutProgramBatchStatusNegative")
$srcSheet.Copy($null, $afterSheet)

# The copy gets auto-named "PutUser (2)" - rename it.
$ws = $wb.Worksheets.Item("PutProgramBatchStatusNegative").Next
$ws.Name = "PutUserIdNegative"

# ---------------------------------------------------------------------
# 2. Drop the columns that aren't part of this sheet: userLogin.loginStatus
#    (N), userLogin.password (O), userLogin.status (P) and, further right,
#    userLogin.roleIds (R), userRoleMaps.userRoleStatus (S),
#    userRoleMaps.roleId (T). Delete right-to-left so the earlier delete
#    doesn't shift the later range.
# ---------------------------------------------------------------------
$ws.Range("R1:T2").EntireColumn.Delete() | Out-Null
$ws.Range("N1:P2").EntireColumn.Delete() | Out-Null

# After the deletes: A..M unchanged, N = userLogin.userLoginEmail,
# O = expectedResponseCode.

# ---------------------------------------------------------------------
# 3. Row 2: turn the copied "PutUser" sample row into the first negative
#    scenario (invalid userId).
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Update User- Invalid userId"
$ws.Range("B2").Value = "U2361111"
$ws.Range("O2").Value = "'404"

# ---------------------------------------------------------------------
# 4. Row 3: duplicate row 2's formatting/values, then adjust for the
#    second negative scenario (invalid timezone).
# ---------------------------------------------------------------------
$ws.Range("A2:O2").Copy($ws.Range("A3:O3"))
$ws.Rows(3).RowHeight = 73

$ws.Range("A3").Value = "Update User- Invalid userId"
$ws.Range("B3").Value = "U236"
$ws.Range("F3").Value = "PSTXX"
$ws.Range("O3").Value = "'400"

# Row 3 didn't inherit the hyperlinks from row 2 - recreate them.
$ws.Hyperlinks.Add($ws.Range("I3"), "https://www.linkedin.com/in/meena/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("N3"), "mailto:Ninjaemailxxx@gmail.com") | Out-Null

# ---------------------------------------------------------------------
# 5. Row 1 (header) is now shorter (fewer columns) so it needs less
#    height than the 21-column "PutUser" header did.
# ---------------------------------------------------------------------
$ws.Rows(1).RowHeight = 29.5

# ---------------------------------------------------------------------
# 6. Column widths that were customized on the new sheet.
# ---------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 31.75
$ws.Columns(14).ColumnWidth = 22.65
$ws.Columns(15).ColumnWidth = 19.25

# ---------------------------------------------------------------------
# 7. View state: scroll/select cell N10 and make this the active tab
#    (this also clears tabSelected on whichever sheet had it before).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("N10").Select()

# ---------------------------------------------------------------------
# 8. "PutUser" sheet selection changed to select the whole used range.
# ---------------------------------------------------------------------
$wsPutUser = $wb.Worksheets.Item("PutUser")
$wsPutUser.Activate()
$wsPutUser.Range("A1:XFD2").Select()

# "PutProgramBatchStatus" sheet is no longer the tab-selected sheet.
$wsPutProgBatch = $wb.Worksheets.Item("PutProgramBatchStatus")
$wsPutProgBatch.Activate()
$wsPutProgBatch.Range("L5").Select()

# Re-activate the new sheet last so it ends up as the active tab.
$ws.Activate()
$ws.Range("N10").Select()
